$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting existing rows 15-39 down to 16-40.
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with the new weekly data point
# (same market/category metadata as the surrounding rows, new date & price figures).
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = "Vega Monumental Concepción"
$ws.Range("C15").Value = "Bíobío"
$ws.Range("D15").Value = 44868
$ws.Range("D15").NumberFormat = $ws.Range("D16").NumberFormat
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 100112026
$ws.Range("G15").Value = "Haba"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 8500
$ws.Range("M15").Value = 8250
$ws.Range("N15").Value = "$/saco 25 kilos"
$ws.Range("O15").Value = "Región de O'Higgins"
$ws.Range("P15").Value = 330
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"
